$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Update Sprint No. from 1 to 2 (Table 1, Row 2, Column 4)
$sprintCell = $tbl.Cell(2, 4)
$sprintCell.Range.Text = "2"

# Update Review Date from 02/09/18 to 02/21/18 (Table 1, Row 3, Column 2)
$dateCell = $tbl.Cell(3, 2)
$dateCell.Range.Text = "02/21/18"
